$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "ADXL375BCCZ-RL7"
$ws.Range("B5").Value = "High G Acc"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 16.4
$ws.Range("F5").Value = "YES"
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = "C481898"
$ws.Range("I5").Value = 5.421

$ws.Range("A8").Select()
